# File Area: Fix support for .xlsx files
#
# The worksheet previously only used B2:D5. This fills in the surrounding
# column A and row 1 cells (blank, but formatted like the rest of the
# table) so the used range grows to A1:D5, gives A5 a real value (5), and
# moves the active selection to B5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C3 already carries the centered / Times New Roman formatting shared by
# the rest of the table's cells - copy it (format only) onto the new
# cells instead of rebuilding the format from scratch, so the new cells
# end up sharing the exact same cell style as their neighbours.
$ws.Range("C3").Copy()

# Row 1: A1:D1 - new blank, styled cells.
$ws.Range("A1:D1").PasteSpecial(-4122)   # xlPasteFormats

# Column A: A2:A4 - new blank, styled cells.
$ws.Range("A2:A4").PasteSpecial(-4122)   # xlPasteFormats

# A5 - styled like the rest, and carries a real numeric value (5).
$ws.Range("A5").PasteSpecial(-4122)      # xlPasteFormats
$ws.Range("A5").Value = 5

$excel.CutCopyMode = $false

# Move the selection to B5 (single cell), matching the new state.
$ws.Range("B5").Select() | Out-Null

Write-Host "Applied File Area .xlsx fix: populated column A / row 1 and re-selected B5"
